# "Task ready for kid piloting"
#
# Adds a new "correctKeyPressNum" column (F) and switches the existing
# "correctKey" (D) / "correctKeyPress" (E) / "allowableKey" (B) values
# from plain digit strings/numbers to "num_*" token strings.
#
# NOTE: a leading single quote in a value that is assigned to a Range is
# interpreted the same way the real Excel UI interprets it: it is consumed
# as a "treat as text" quote-prefix marker (and the cell picks up Excel's
# built in quotePrefix cell style) rather than becoming part of the stored
# text. To store a value that itself starts with a literal apostrophe
# (e.g. 'num_7') two leading quotes are used, which collapse to one
# literal leading apostrophe in the stored text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("correctKey"): ['1'] / ['4'] / ['7'] / ['0']  ->  ['num_7'] / ['num_4'] / ['num_1'] / ['num_0']
$ws.Range("D2").Value = "'['num_7']"
$ws.Range("D3").Value = "'['num_4']"
$ws.Range("D4").Value = "'['num_1']"
$ws.Range("D5").Value = "'['num_0']"

# Column B ("allowableKey"): same list value on every data row.
$ws.Range("B2:B5").Value = "'['num_7', 'num_4', 'num_1','num_0']"

# New column F ("correctKeyPressNum"): plain num_* tokens, default style.
$ws.Range("F1").Value = "correctKeyPressNum"
$ws.Range("F2").Value = "num_7"
$ws.Range("F3").Value = "num_4"
$ws.Range("F4").Value = "num_1"
$ws.Range("F5").Value = "num_0"

# Column E ("correctKeyPress"): was bare numbers 7/4/1/0, now quoted tokens.
$ws.Range("E2").Value = "''num_7'"
$ws.Range("E3").Value = "''num_4'"
$ws.Range("E4").Value = "''num_1'"
$ws.Range("E5").Value = "''num_0'"

# Move the sheet's selection down to below the new table, one column to
# the right of where it used to be (B6 -> E6).
$null = $ws.Range("E6").Select()
